$wb = $excel.ActiveWorkbook

# --- Sheet "Events": fix header A1 (was garbled "jjjjjj") ---
$wsEvents = $wb.Worksheets.Item("Events")
$wsEvents.Range("A1").Value = "start"
$wsEvents.Activate()
$wsEvents.Range("A2").Select()

# --- Sheet "Artworks": insert a new "end" column after "start" (year), ---
# --- and rename the "artwork" column to "name" to match the common schema ---
$wsArtworks = $wb.Worksheets.Item("Artworks")
$wsArtworks.Columns.Item(2).Insert()
$wsArtworks.Range("A1").Value = "start"
$wsArtworks.Range("B1").Value = "end"
$wsArtworks.Range("E1").Value = "name"
$wsArtworks.Activate()
$wsArtworks.Range("B3").Select()

# --- Sheet "Styles": rename the "style" column to "name" ---
$wsStyles = $wb.Worksheets.Item("Styles")
$wsStyles.Range("C1").Value = "name"
$wsStyles.Activate()
$wsStyles.Range("C1").Select()

$wsEvents.Activate()
